$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume(1h) (E) columns with refreshed crypto data ---
# D cells hold text-formatted numbers (e.g. "9.490", "0.06680"); force text
# format before assignment so Excel does not coerce them into real numbers
# and strip meaningful trailing/leading zeros.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.352.34"
$ws.Range("E2").Value = "  +0.10%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.565.33"
$ws.Range("E3").Value = "  +0.11%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.70%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.007"
$ws.Range("E5").Value = "  +0.48%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "288.76"
$ws.Range("E6").Value = "  -0.38%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3726"
$ws.Range("E7").Value = "  -0.36%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.21"
$ws.Range("E8").Value = "  -0.21%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3348"
$ws.Range("E9").Value = "  -2.27%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07378"
$ws.Range("E10").Value = "  -3.45%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.111"
$ws.Range("E11").Value = "  -4.64%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.008"
$ws.Range("E12").Value = "  +0.72%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.61"
$ws.Range("E13").Value = "  -3.83%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.825"
$ws.Range("E14").Value = "  -3.31%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.798"
$ws.Range("E15").Value = "  -1.92%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.569.28"
$ws.Range("E16").Value = "  +0.30%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001097"
$ws.Range("E17").Value = "  -2.62%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "88.88"
$ws.Range("E18").Value = "  -1.03%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06680"
$ws.Range("E19").Value = "  -0.68%  "

$ws.Range("E20").Value = "  +0.60%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.111"
$ws.Range("E21").Value = "  -1.95%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.97"
$ws.Range("E22").Value = "  -3.55%  "

$ws.Range("E23").Value = "  -1.65%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.350.18"
$ws.Range("E24").Value = "  +0.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.369"
$ws.Range("E25").Value = "  -1.30%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.474"
$ws.Range("E26").Value = "  -11.66%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.73"
$ws.Range("E27").Value = "  -2.50%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "146.12"
$ws.Range("E28").Value = "  +0.21%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.982"
$ws.Range("E29").Value = "  -0.12%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.72"
$ws.Range("E30").Value = "  -1.35%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.742.68"
$ws.Range("E31").Value = "  +0.36%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.984"
$ws.Range("E32").Value = "  -1.48%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9677"
$ws.Range("E33").Value = "  -4.14%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.835"
$ws.Range("E34").Value = "  -5.39%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08362"
$ws.Range("E35").Value = "  -1.80%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.490"
$ws.Range("E36").Value = "  -5.73%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.381"
$ws.Range("E37").Value = "  +5.72%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02431"
$ws.Range("E38").Value = "  -4.40%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2242"
$ws.Range("E39").Value = "  -3.10%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06297"
$ws.Range("E40").Value = "  -1.53%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.303"
$ws.Range("E41").Value = "  -3.41%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6097"
$ws.Range("E42").Value = "  -3.87%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.007"
$ws.Range("E43").Value = "  +0.60%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.83"
$ws.Range("E44").Value = "  -7.50%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.73"
$ws.Range("E45").Value = "  -1.95%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.778"
$ws.Range("E46").Value = "  +0.73%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5687"
$ws.Range("E47").Value = "  -4.93%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.211"
$ws.Range("E50").Value = "  -4.09%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07287"
$ws.Range("E51").Value = "  +0.52%  "

# --- Rows 48/49: NEARProtocol and Quant swap list positions, each with fresh data ---
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.03"
$ws.Range("E48").Value = "  +0.37%  "

$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.009"
$ws.Range("E49").Value = "  -3.85%  "
